$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete now-removed rows (previously rows 11-13, Resolving-Mac -> Resolving-Mac target)
$ws.Rows("11:13").Delete()

# Row 2: ECs -> ECs
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Pdgfb"
$ws.Range("C2").Value = "Art1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 45.01222466666667
$ws.Range("H2").Value = 135.036674
$ws.Range("I2").Value = 0.7482903203664146
$ws.Range("J2").Value = 0.7482903203664146
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.2283436666666667
$ws.Range("N2").Value = 0.6850310000000001
$ws.Range("O2").Value = 0.09370018555900643
$ws.Range("P2").Value = 0.09370018555900644
$ws.Range("Q2").Value = 10.27825642521045
$ws.Range("R2").Value = 92.50430782689401
$ws.Range("S2").Value = 0.07011494187034142
$ws.Range("T2").Value = 0.07011494187034142

# Row 3: ECs -> FAPs
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Pdgfb"
$ws.Range("C3").Value = "Art1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 45.01222466666667
$ws.Range("H3").Value = 135.036674
$ws.Range("I3").Value = 0.7482903203664146
$ws.Range("J3").Value = 0.7482903203664146
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.5088386666666667
$ws.Range("N3").Value = 1.526516
$ws.Range("O3").Value = 0.2088005250255714
$ws.Range("P3").Value = 0.2088005250255715
$ws.Range("Q3").Value = 22.90396038308711
$ws.Range("R3").Value = 206.135643447784
$ws.Range("S3").Value = 0.1562434117640604
$ws.Range("T3").Value = 0.1562434117640605

# Row 4: ECs -> MuSCs
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Pdgfb"
$ws.Range("C4").Value = "Art1"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 45.01222466666667
$ws.Range("H4").Value = 135.036674
$ws.Range("I4").Value = 0.7482903203664146
$ws.Range("J4").Value = 0.7482903203664146
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 1.699778333333333
$ws.Range("N4").Value = 5.099335
$ws.Range("O4").Value = 0.6974992894154221
$ws.Range("P4").Value = 0.6974992894154222
$ws.Range("Q4").Value = 76.51080422353223
$ws.Range("R4").Value = 688.5972380117901
$ws.Range("S4").Value = 0.5219319667320127
$ws.Range("T4").Value = 0.5219319667320128

# Row 5: MuSCs -> ECs
$ws.Range("A5").Value = "MuSCs"
$ws.Range("B5").Value = "Pdgfb"
$ws.Range("C5").Value = "Art1"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 2.766295666666667
$ws.Range("H5").Value = 8.298887000000001
$ws.Range("I5").Value = 0.04598733535094824
$ws.Range("J5").Value = 0.04598733535094825
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.2283436666666667
$ws.Range("N5").Value = 0.6850310000000001
$ws.Range("O5").Value = 0.09370018555900643
$ws.Range("P5").Value = 0.09370018555900644
$ws.Range("Q5").Value = 0.6316660956107779
$ws.Range("R5").Value = 5.684994860497
$ws.Range("S5").Value = 0.004309021855748106
$ws.Range("T5").Value = 0.004309021855748107

# Row 6: MuSCs -> FAPs
$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Pdgfb"
$ws.Range("C6").Value = "Art1"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 2.766295666666667
$ws.Range("H6").Value = 8.298887000000001
$ws.Range("I6").Value = 0.04598733535094824
$ws.Range("J6").Value = 0.04598733535094825
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.5088386666666667
$ws.Range("N6").Value = 1.526516
$ws.Range("O6").Value = 0.2088005250255714
$ws.Range("P6").Value = 0.2088005250255715
$ws.Range("Q6").Value = 1.407598198632444
$ws.Range("R6").Value = 12.668383787692
$ws.Range("S6").Value = 0.009602179765805014
$ws.Range("T6").Value = 0.009602179765805017

# Row 7: MuSCs -> MuSCs
$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Pdgfb"
$ws.Range("C7").Value = "Art1"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 2.766295666666667
$ws.Range("H7").Value = 8.298887000000001
$ws.Range("I7").Value = 0.04598733535094824
$ws.Range("J7").Value = 0.04598733535094825
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 1.699778333333333
$ws.Range("N7").Value = 5.099335
$ws.Range("O7").Value = 0.6974992894154221
$ws.Range("P7").Value = 0.6974992894154222
$ws.Range("Q7").Value = 4.702089437793889
$ws.Range("R7").Value = 42.318804940145
$ws.Range("S7").Value = 0.03207613372939511
$ws.Range("T7").Value = 0.03207613372939513

# Row 8: Resolving-Mac -> ECs
$ws.Range("A8").Value = "Resolving-Mac"
$ws.Range("B8").Value = "Pdgfb"
$ws.Range("C8").Value = "Art1"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 12.37490333333333
$ws.Range("H8").Value = 37.12471
$ws.Range("I8").Value = 0.2057223442826371
$ws.Range("J8").Value = 0.2057223442826371
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.2283436666666667
$ws.Range("N8").Value = 0.6850310000000001
$ws.Range("O8").Value = 0.09370018555900643
$ws.Range("P8").Value = 0.09370018555900644
$ws.Range("Q8").Value = 2.82573080177889
$ws.Range("R8").Value = 25.43157721601
$ws.Range("S8").Value = 0.01927622183291691
$ws.Range("T8").Value = 0.01927622183291691

# Row 9: Resolving-Mac -> FAPs
$ws.Range("A9").Value = "Resolving-Mac"
$ws.Range("B9").Value = "Pdgfb"
$ws.Range("C9").Value = "Art1"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 12.37490333333333
$ws.Range("H9").Value = 37.12471
$ws.Range("I9").Value = 0.2057223442826371
$ws.Range("J9").Value = 0.2057223442826371
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.5088386666666667
$ws.Range("N9").Value = 1.526516
$ws.Range("O9").Value = 0.2088005250255714
$ws.Range("P9").Value = 0.2088005250255715
$ws.Range("Q9").Value = 6.296829312262222
$ws.Range("R9").Value = 56.67146381036
$ws.Range("S9").Value = 0.042954933495706
$ws.Range("T9").Value = 0.042954933495706

# Row 10: Resolving-Mac -> MuSCs
$ws.Range("A10").Value = "Resolving-Mac"
$ws.Range("B10").Value = "Pdgfb"
$ws.Range("C10").Value = "Art1"
$ws.Range("D10").Value = "MuSCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 12.37490333333333
$ws.Range("H10").Value = 37.12471
$ws.Range("I10").Value = 0.2057223442826371
$ws.Range("J10").Value = 0.2057223442826371
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 1.699778333333333
$ws.Range("N10").Value = 5.099335
$ws.Range("O10").Value = 0.6974992894154221
$ws.Range("P10").Value = 0.6974992894154222
$ws.Range("Q10").Value = 21.03459256309445
$ws.Range("R10").Value = 189.31133306785
$ws.Range("S10").Value = 0.1434911889540142
$ws.Range("T10").Value = 0.1434911889540142
